# Update Valmir's weekly schedule sheet: move several class blocks to the
# days/slots reflected in the latest version ("atualizado" per commit msg).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '[-, ''MCT-1A-Metrologia'', -, -]'
$ws.Range("C3").Value = '[-, ''MCT-1A-Metrologia'', -, -]'
$ws.Range("F3").Value = '[-, -, -, ''MEC-3B-Calderaria'']'
$ws.Range("C4").Value = '[-, ''MCT-1A-Metrologia'', -, -]'
$ws.Range("F4").Value = '[-, -, -, Victor-Usin. CNC-3B]'
$ws.Range("C6").Value = '-'
$ws.Range("F6").Value = '[-, -, -, ''MEC-3B-Calderaria'']'
$ws.Range("C7").Value = '-'
$ws.Range("F7").Value = '[-, -, -, ''MEC-3B-Calderaria'']'
$ws.Range("D8").Value = '[''MCT-1A-Metrologia'', -, -, -]'
$ws.Range("B11").Value = '[-, -, ''MEC-3A-Calderaria'', -]'
$ws.Range("E11").Value = '-'
$ws.Range("B12").Value = '[-, -, ''MEC-3A-Calderaria'', -]'
$ws.Range("D12").Value = '[-, -, -, -]'
$ws.Range("E12").Value = '-'
$ws.Range("B14").Value = '[-, -, ''MEC-3A-Calderaria'', -]'
$ws.Range("D14").Value = '[-, -, -, -]'
$ws.Range("E14").Value = '-'
$ws.Range("B15").Value = '[-, -, ''MEC-3A-Calderaria'', -]'
$ws.Range("E15").Value = '-'
$ws.Range("B18").Value = '[-, ''MEC-1NB-Caldeiraria'', ''MEC-1NA-Tec. Mat. Não Metal.'', -]'
$ws.Range("C18").Value = '-'
$ws.Range("F18").Value = '-'
$ws.Range("B19").Value = '[-, ''MEC-1NB-Caldeiraria'', ''MEC-1NA-Tec. Mat. Não Metal.'', -]'
$ws.Range("C19").Value = '-'
$ws.Range("F19").Value = '-'
$ws.Range("B20").Value = '[-, ''MEC-1NB-Caldeiraria'', -, -]'
$ws.Range("C20").Value = '-'
$ws.Range("D20").Value = '[''MEC-1NA-Tec. Mat. Não Metal.'', -, -, -]'
$ws.Range("F20").Value = '-'
$ws.Range("B21").Value = '[-, ''MEC-1NB-Caldeiraria'', -, -]'
$ws.Range("C21").Value = '-'
$ws.Range("D21").Value = '[''MEC-1NA-Tec. Mat. Não Metal.'', -, -, -]'
$ws.Range("F21").Value = '-'
